# Updated cryptos list on Thu May 16 19:30:44 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.492.66"
$ws.Range("E2").Value = "  -0.72%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.953.66"
$ws.Range("E3").Value = "  -1.92%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.13%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.98%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.45%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.12%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.44%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.946.21"
$ws.Range("E9").Value = "  -2.02%  "

# Row 10 - Toncoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.54%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -1.50%  "

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.16%  "

# Row 13 - ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.86%  "

# Row 14 - Avalanche
$ws.Range("E14").Value = "  +0.50%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  -0.74%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "65.408.59"
$ws.Range("E16").Value = "  -0.70%  "

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "3.443.12"
$ws.Range("E17").Value = "  -1.78%  "

# Row 18 - Polkadot
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.61%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "2.952.10"
$ws.Range("E19").Value = "  -1.93%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.46%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "447.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.47%  "

# Row 22 - Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.692"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.23%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  -1.11%  "

# Row 24 - Litecoin
$ws.Range("E24").Value = "  +0.88%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  -1.13%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.10%  "

# Row 27 - RenderToken
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.08%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  +0.02%  "

# Row 29 - NEARProtocol
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.42%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  -0.19%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.81%  "

# Row 32 - PEPE
$ws.Range("E32").Value = "  -2.15%  "

# Row 33 - EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.61%  "

# Row 34 - Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.48%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  +0.18%  "

# Row 36 - Mantle
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.978"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.53%  "

# Row 37 - Filecoin
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.84%  "

# Row 38 - OKB
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.46%  "

# Row 39 - was Arweave, now Stacks (row content swapped with row 40, with updated figures)
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.10%  "

# Row 40 - was Stacks, now Arweave (row content swapped with row 39, with updated figures)
$ws.Range("B40").Value = "Arweave"
$ws.Range("C40").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "44.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.84%  "

# Row 41 - dogwifhat
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.78%  "

# Row 42 - Kaspa
$ws.Range("E42").Value = "  -1.07%  "

# Row 43 - TheGraph
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.299"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "

# Row 44 - Cosmos
$ws.Range("E44").Value = "  +0.30%  "

# Row 45 - Bittensor
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "386.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.37%  "

# Row 46 - VeChain
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0353"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "

# Row 47 - Maker
$ws.Range("D47").Value = "2.715.82"
$ws.Range("E47").Value = "  -1.87%  "

# Row 48 - Monero
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.73%  "

# Row 50 - ThetaToken
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.40%  "

# Row 51 - was InjectiveProtocol, now Stellar
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.106"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "

